$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.784.57'
$ws.Range('E2').Value = '  +0.03%  '
$ws.Range('D3').Value = '1.648.95'
$ws.Range('E3').Value = '  +0.00%  '
$ws.Range('E4').Value = '  +0.71%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '216.90'
$ws.Range('D5').NumberFormat = 'General'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.81%  '
$ws.Range('E6').Value = '  +0.22%  '
$ws.Range('E7').Value = '  +0.68%  '
$ws.Range('E8').Value = '  +0.20%  '
$ws.Range('E9').Value = '  -0.31%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.31'
$ws.Range('D10').NumberFormat = 'General'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.02%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0846'
$ws.Range('D11').NumberFormat = 'General'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.13%  '
$ws.Range('D12').Value = '1.873.36'
$ws.Range('E12').Value = '  -0.26%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.655.01'
$ws.Range('E13').Value = '  +0.41%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.22'
$ws.Range('D14').NumberFormat = 'General'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.38%  '
$ws.Range('E15').Value = '  +0.20%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.78'
$ws.Range('D16').NumberFormat = 'General'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.74%  '
$ws.Range('D17').Value = '26.795.62'
$ws.Range('E17').Value = '  -0.08%  '
$ws.Range('D18').Value = '0.0₃0746'
$ws.Range('E18').Value = '  -0.27%  '
$ws.Range('E19').Value = '  -0.65%  '
$ws.Range('E20').Value = '  +0.62%  '
$ws.Range('E21').Value = '  +0.34%  '
$ws.Range('E22').Value = '  +16.88%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.33'
$ws.Range('D23').NumberFormat = 'General'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.64%  '
$ws.Range('E24').Value = '  +0.04%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '145.75'
$ws.Range('D25').NumberFormat = 'General'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.39%  '
$ws.Range('E26').Value = '  +0.70%  '
$ws.Range('E27').Value = '  -0.85%  '
$ws.Range('E28').Value = '  +3.85%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.83'
$ws.Range('D29').NumberFormat = 'General'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.14%  '
$ws.Range('E30').Value = '  +0.56%  '
$ws.Range('E31').Value = '  +0.66%  '
$ws.Range('E32').Value = '  -0.93%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.03'
$ws.Range('D33').NumberFormat = 'General'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.43%  '
$ws.Range('B34').Value = 'LidoDAOToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.55'
$ws.Range('D34').NumberFormat = 'General'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.89%  '
$ws.Range('B35').Value = 'Maker'
$ws.Range('C35').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D35').Value = '1.277.96'
$ws.Range('E35').Value = '  -0.15%  '
$ws.Range('E36').Value = '  +1.99%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0179'
$ws.Range('D37').NumberFormat = 'General'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.20%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.543'
$ws.Range('D38').NumberFormat = 'General'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +5.18%  '
$ws.Range('E39').Value = '  +2.82%  '
$ws.Range('E40').Value = '  +0.72%  '
$ws.Range('E41').Value = '  +1.53%  '
$ws.Range('E42').Value = '  -1.38%  '
$ws.Range('E43').Value = '  +1.02%  '
$ws.Range('D44').Value = '1.798.48'
$ws.Range('E44').Value = '  +0.62%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '92.20'
$ws.Range('D45').NumberFormat = 'General'
$ws.Range('D45').Style = 'Normal'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '59.61'
$ws.Range('D46').NumberFormat = 'General'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +6.42%  '
$ws.Range('E47').Value = '  +1.43%  '
$ws.Range('E48').Value = '  +1.20%  '
$ws.Range('E49').Value = '  +0.31%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.82'
$ws.Range('D50').NumberFormat = 'General'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.63%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0984'
$ws.Range('D51').NumberFormat = 'General'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.59%  '
